# Scheduled-runner style data refresh: update the FFXIV leve market-price /
# profit columns (H..N) on a handful of rows across the eight item-crafting
# sheets with freshly pulled values.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 10166.667
$ws.Range("J16").Value = 11500
$ws.Range("L16").Value = 11500
$ws.Range("N16").Value = -11960
$ws.Range("H21").Value = 7008.5
$ws.Range("I21").Value = 1017
$ws.Range("K21").Value = 1017
$ws.Range("M21").Value = -549
$ws.Range("H23").Value = 7008.5
$ws.Range("I23").Value = 1017
$ws.Range("K23").Value = 1017
$ws.Range("M23").Value = -783
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("N26").Value = 0
$ws.Range("H33").Value = 1321.5333
$ws.Range("I33").Value = 126.75
$ws.Range("K33").Value = 126.75
$ws.Range("M33").Value = 102.25
$ws.Range("H40").Value = 1689.5555
$ws.Range("I40").Value = 1615
$ws.Range("K40").Value = 1615
$ws.Range("M40").Value = -1440
$ws.Range("H43").Value = 1350
$ws.Range("I43").Value = 1380
$ws.Range("J43").Value = 1200
$ws.Range("K43").Value = 1380
$ws.Range("L43").Value = 1200
$ws.Range("M43").Value = -1311
$ws.Range("N43").Value = -1338
$ws.Range("H106").Value = 1618.5333
$ws.Range("I106").Value = 1670.5714
$ws.Range("K106").Value = 1670.5714
$ws.Range("M106").Value = -1039.5714
$ws.Range("H116").Value = 1917.9231
$ws.Range("I116").Value = 1179.4546
$ws.Range("J116").Value = 5979.5
$ws.Range("K116").Value = 1179.4546
$ws.Range("L116").Value = 5979.5
$ws.Range("M116").Value = 2262.5454
$ws.Range("N116").Value = -12863.5
$ws.Range("H132").Value = 7357528.5
$ws.Range("I132").Value = 7817253
$ws.Range("J132").Value = 1939
$ws.Range("K132").Value = 23451759
$ws.Range("L132").Value = 5817
$ws.Range("M132").Value = -23449229
$ws.Range("N132").Value = -10877
$ws.Range("L26").Value = ""   # cell removed from row 26 in the diff

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 38069.17
$ws.Range("I32").Value = 7408.6665
$ws.Range("J32").Value = 92175.94
$ws.Range("K32").Value = 7408.6665
$ws.Range("L32").Value = 92175.94
$ws.Range("M32").Value = -7121.6665
$ws.Range("N32").Value = -92749.94
$ws.Range("H110").Value = 25056110
$ws.Range("I110").Value = 29477312
$ws.Range("J110").Value = 2637
$ws.Range("K110").Value = 29477312
$ws.Range("L110").Value = 2637
$ws.Range("M110").Value = -29475267
$ws.Range("N110").Value = -6727

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16921.445
$ws.Range("I31").Value = 33341.773
$ws.Range("J31").Value = 1949.9706
$ws.Range("K31").Value = 33341.773
$ws.Range("L31").Value = 1949.9706
$ws.Range("M31").Value = -33046.773
$ws.Range("N31").Value = -2539.9706
$ws.Range("H34").Value = 16921.445
$ws.Range("I34").Value = 33341.773
$ws.Range("J34").Value = 1949.9706
$ws.Range("K34").Value = 33341.773
$ws.Range("L34").Value = 1949.9706
$ws.Range("M34").Value = -33139.773
$ws.Range("N34").Value = -2353.9706
$ws.Range("H80").Value = 15752
$ws.Range("J80").Value = 15752
$ws.Range("L80").Value = 15752
$ws.Range("N80").Value = -17998
$ws.Range("H83").Value = 15752
$ws.Range("J83").Value = 15752
$ws.Range("L83").Value = 47256
$ws.Range("N83").Value = -58488
$ws.Range("H134").Value = 1765.4
$ws.Range("I134").Value = 1360.6666
$ws.Range("J134").Value = 3384.3333
$ws.Range("K134").Value = 4081.9998
$ws.Range("L134").Value = 10152.9999
$ws.Range("M134").Value = -1546.9998
$ws.Range("N134").Value = -15222.9999

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 2666.5833
$ws.Range("J33").Value = 2400
$ws.Range("L33").Value = 14400
$ws.Range("N33").Value = -14966
$ws.Range("H81").Value = 33335858
$ws.Range("I81").Value = 993.3333
$ws.Range("J81").Value = 47622228
$ws.Range("K81").Value = 2979.9999
$ws.Range("L81").Value = 142866684
$ws.Range("M81").Value = -1856.9999
$ws.Range("N81").Value = -142868930
$ws.Range("H84").Value = 33335858
$ws.Range("I84").Value = 993.3333
$ws.Range("J84").Value = 47622228
$ws.Range("K84").Value = 8939.9997
$ws.Range("L84").Value = 428600052
$ws.Range("M84").Value = -3323.9997
$ws.Range("N84").Value = -428611284
$ws.Range("H107").Value = 1283.4193
$ws.Range("I107").Value = 1001
$ws.Range("J107").Value = 1337.7307
$ws.Range("K107").Value = 3003
$ws.Range("L107").Value = 4013.1921
$ws.Range("M107").Value = -1083
$ws.Range("N107").Value = -7853.1921

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("M18").Value = -207
$ws.Range("H132").Value = 3456.75
$ws.Range("I132").Value = 3096.2942
$ws.Range("J132").Value = 5499.3335
$ws.Range("K132").Value = 9288.882599999999
$ws.Range("L132").Value = 16498.0005
$ws.Range("M132").Value = -6758.882599999999
$ws.Range("N132").Value = -21558.0005

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("H40").Value = 80484.62
$ws.Range("I40").Value = 94163.73
$ws.Range("K40").Value = 94163.73
$ws.Range("M40").Value = -94027.73
$ws.Range("H122").Value = 3211.889
$ws.Range("I122").Value = 3211.889
$ws.Range("K122").Value = 9635.667000000001
$ws.Range("M122").Value = -7185.667000000001
$ws.Range("H127").Value = 41312.5
$ws.Range("J127").Value = 41312.5
$ws.Range("L127").Value = 41312.5
$ws.Range("N127").Value = -51232.5
$ws.Range("M30").Value = ""   # cell removed from row 30 in the diff

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 986.3333
$ws.Range("I122").Value = 979.5
$ws.Range("K122").Value = 2938.5
$ws.Range("M122").Value = -488.5
$ws.Range("H132").Value = 4509.3687
$ws.Range("I132").Value = 7545.25
$ws.Range("J132").Value = 2301.4546
$ws.Range("K132").Value = 22635.75
$ws.Range("L132").Value = 6904.3638
$ws.Range("M132").Value = -20105.75
$ws.Range("N132").Value = -11964.3638
$ws.Range("H136").Value = 647.25
$ws.Range("I136").Value = 452.51724
$ws.Range("J136").Value = 1023.73334
$ws.Range("K136").Value = 1357.55172
$ws.Range("L136").Value = 3071.20002
$ws.Range("M136").Value = 1192.44828
$ws.Range("N136").Value = -8171.20002
